$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("FAPs", "Cort", "Ghsr", "ECs", 2, 0.6666666666666666, 0.2586766666666667, 0.77603, 0.9158046909573684, 0.9158046909573684, 1, 0.3333333333333333, 0.2864683333333333, 0.859405, 0.1191104989459046, 0.1191104989459046, 0.07410267357222222, 0.66692406215, 0.1090819536769321, 0.1090819536769321)
    ,@("FAPs", "Cort", "Ghsr", "FAPs", 2, 0.6666666666666666, 0.2586766666666667, 0.77603, 0.9158046909573684, 0.9158046909573684, 3, 1, 1.671010333333333, 5.013031, 0.6947883985330395, 0.6947883985330396, 0.4322513829922222, 3.89026244693, 0.6362904745993152, 0.6362904745993153)
    ,@("FAPs", "Cort", "Ghsr", "Inflammatory-Mac", 2, 0.6666666666666666, 0.2586766666666667, 0.77603, 0.9158046909573684, 0.9158046909573684, 1, 0.3333333333333333, 0.07779633333333334, 0.233389, 0.03234689144057309, 0.03234689144057309, 0.02012409618555556, 0.18111686567, 0.02962343491916559, 0.02962343491916559)
    ,@("FAPs", "Cort", "Ghsr", "MuSCs", 2, 0.6666666666666666, 0.2586766666666667, 0.77603, 0.9158046909573684, 0.9158046909573684, 3, 1, 0.3697886666666667, 1.109366, 0.1537542110804829, 0.1537542110804829, 0.09565569966444445, 0.8609012969800001, 0.1408088277619556, 0.1408088277619556)
    ,@("Resolving-Mac", "Cort", "Ghsr", "ECs", 1, 0.3333333333333333, 0.02378166666666667, 0.07134500000000001, 0.08419530904263167, 0.08419530904263166, 1, 0.3333333333333333, 0.2864683333333333, 0.859405, 0.1191104989459046, 0.1191104989459046, 0.00681269441388889, 0.06131424972500001, 0.01002854526897249, 0.01002854526897249)
    ,@("Resolving-Mac", "Cort", "Ghsr", "FAPs", 1, 0.3333333333333333, 0.02378166666666667, 0.07134500000000001, 0.08419530904263167, 0.08419530904263166, 3, 1, 1.671010333333333, 5.013031, 0.6947883985330395, 0.6947883985330396, 0.03973941074388889, 0.357654696695, 0.0584979239337244, 0.0584979239337244)
    ,@("Resolving-Mac", "Cort", "Ghsr", "Inflammatory-Mac", 1, 0.3333333333333333, 0.02378166666666667, 0.07134500000000001, 0.08419530904263167, 0.08419530904263166, 1, 0.3333333333333333, 0.07779633333333334, 0.233389, 0.03234689144057309, 0.03234689144057309, 0.001850126467222223, 0.016651138205, 0.002723456521407509, 0.002723456521407508)
    ,@("Resolving-Mac", "Cort", "Ghsr", "MuSCs", 1, 0.3333333333333333, 0.02378166666666667, 0.07134500000000001, 0.08419530904263167, 0.08419530904263166, 3, 1, 0.3697886666666667, 1.109366, 0.1537542110804829, 0.1537542110804829, 0.008794190807777781, 0.07914771727000001, 0.01294538331852728, 0.01294538331852728)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
